# Add two new tasks to the TodoList worksheet (rows 48 and 49).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 48: Init Exceptions ------------------------------------------------
$ws.Range("A48").Value = "Init Exceptions"
$ws.Range("B48").Value = "Design Issue"
$ws.Range("C48").Value = "We should have a new exception called InitializationException that is thrown whenever the init function doesn't succeed. It collects the actual reason or exception why it didn't succeed and puts that as the message of the exception."
$ws.Range("D48").Value = 40948
$ws.Range("E48").Value = "Yael"
$ws.Range("F48").Value = "High"

# --- Row 49: Dlog group-- isIdentity() --------------------------------------
$ws.Range("A49").Value = "Dlog group-- isIdentity()"
$ws.Range("B49").Value = "New Feature"
$ws.Range("C49").Value = "Add  a new function to DlogGroup interface: boolean isIdentity(GroupElement gEl)"
$ws.Range("D49").Value = 40949
$ws.Range("E49").Value = "Yael"
$ws.Range("F49").Value = "High"

# Mirror the style/formatting of the existing alternating task blocks
# (row 45 = orange block, row 46/47 = green block), so the new
# rows continue the alternation: row48 = orange block, row49 = green block.
$ws.Range("A45:G45").Copy()
$ws.Range("A48:G48").PasteSpecial(-4122)
$ws.Range("A46:G46").Copy()
$ws.Range("A49:G49").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(48).RowHeight = $ws.Rows.Item(45).RowHeight

# Update the frozen pane / selection to match the new bottom of the list.
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("B27").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("A48").Select()
